# "update as per coach review" -- add the coach's Version 4 / Version 5
# review rows to the RTM and Design Document sheets, and leave the
# Design Document sheet as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# RTM sheet
# ---------------------------------------------------------------
$rtm = $wb.Worksheets.Item("RTM")

# Insert a new "Version 1" row right after the header, before the
# existing Version 2 / Version 3 rows (they shift down to rows 3-4).
$rtm.Rows.Item(2).Insert() | Out-Null
$rtm.Range("A2:D2").ClearFormats() | Out-Null

$rtm.Range("A2").Value = "Version 1"
$rtm.Range("C2").Value = "Amr"
$rtm.Rows.Item(2).RowHeight = 45

$rtm.Range("A2").Interior.Color = 16777215
$rtm.Range("C2").Interior.Color = 16777215
$rtm.Range("D2").Interior.Color = 16777215
$rtm.Range("B2").WrapText = $true

# New "Version 4" rows appended at the bottom (rows 5-8)
$rtm.Range("A5").Value = "Version 4"
$rtm.Range("B5").Value = "add classes names in code " + [char]10 + "beside each feature"
$rtm.Range("C5").Value = "Amr"

$rtm.Range("A6").Value = "Version 4"
$rtm.Range("B6").Value = "add use case beside main " + [char]10 + "functions"
$rtm.Range("C6").Value = "Amr"

$rtm.Range("A7").Value = "Version 4"
$rtm.Range("B7").Value = "add erd beside any function " + [char]10 + "related to database."
$rtm.Range("C7").Value = "Amr"

$rtm.Range("A8").Value = "Version 4"
$rtm.Range("B8").Value = "add wire frames"
$rtm.Range("C8").Value = "Amr"

$rtm.Range("B5:B7").WrapText = $true
$rtm.Rows.Item(5).RowHeight = 30
$rtm.Rows.Item(6).RowHeight = 30
$rtm.Rows.Item(7).RowHeight = 30

# ---------------------------------------------------------------
# Design Document sheet
# ---------------------------------------------------------------
$dd = $wb.Worksheets.Item("Design Document")

$dd.Range("A5").Value = "Version 5"
$dd.Range("B5").Value = "sequence diagrams are required to be generic" + [char]10 + "not high level."
$dd.Range("C5").Value = "Amr"

$dd.Range("A6").Value = "Version 5"
$dd.Range("B6").Value = "change IDs of navigation function flowchart" + [char]10 + "from low level to high level"
$dd.Range("C6").Value = "Amr"

$dd.Range("B5:B6").WrapText = $true
$dd.Rows.Item(5).RowHeight = 30
$dd.Rows.Item(6).RowHeight = 30

# ---------------------------------------------------------------
# Make the Design Document tab the active/selected sheet, with B7
# selected (next empty comment row), matching the reviewer's
# last-saved view state.
# ---------------------------------------------------------------
$dd.Activate() | Out-Null
$dd.Range("B7").Select() | Out-Null

Write-Output "done"
